# Generate Report for Archive
#
# Updates the localization "Status" from "Ready for handoff" to
# "In Translation" across the Overview sheet (zh-cn / de-de columns) and
# each per-locale status sheet, then lets the report's Status columns
# shrink to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$wsOverview.Range("E2:F3").Replace($oldStatus, $newStatus) | Out-Null

# --- Per-locale sheets: Status column (col C) ---
$wsZhCn.Range("C2:C3").Replace($oldStatus, $newStatus) | Out-Null
$wsDeDe.Range("C2:C3").Replace($oldStatus, $newStatus) | Out-Null

# --- Re-fit the Status columns now that the text is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
